$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Cfh -> Sell, ECs): ligand avg/total expression were recomputed with
# the new TPM values, cascading into the derived-specificity and edge columns.
$ws.Range("G2").Value = 0.4128076666666667
$ws.Range("H2").Value = 1.238423
$ws.Range("I2").Value = 0.001366259689176221
$ws.Range("J2").Value = 0.001366259689176221
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.002279333333333333
$ws.Range("N2").Value = 0.006838
$ws.Range("Q2").Value = 0.0009409262748888889
$ws.Range("R2").Value = 0.008468336474000001
$ws.Range("S2").Value = 0.001366259689176221
$ws.Range("T2").Value = 0.001366259689176221

# Row 3 (Cfh -> Sell, FAPs)
$ws.Range("I3").Value = 0.90768474543873
$ws.Range("J3").Value = 0.9076847454387301
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.002279333333333333
$ws.Range("N3").Value = 0.006838
$ws.Range("Q3").Value = 0.6251113408857777
$ws.Range("R3").Value = 5.626002067972
$ws.Range("S3").Value = 0.90768474543873
$ws.Range("T3").Value = 0.9076847454387301

# Row 4 (Cfh -> Sell, MuSCs)
$ws.Range("I4").Value = 0.09094899487209368
$ws.Range("J4").Value = 0.09094899487209368
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.002279333333333333
$ws.Range("N4").Value = 0.006838
$ws.Range("Q4").Value = 0.06263545622244444
$ws.Range("R4").Value = 0.563719106002
$ws.Range("S4").Value = 0.09094899487209368
$ws.Range("T4").Value = 0.09094899487209368
